# CE342_A.xlsx — Turma A: add "Resenha Regime de Metas" scores as a new
# column L, shifting the existing "Email" column from L to M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# M1 should end up with the exact same formatting (style) the "Email"
# header currently has in L1, so snapshot that formatting onto M1 first.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# New "Resenha Regime de Metas" scores for rows 2..44 (row -> score).
# $null marks rows that are blank for this assignment (same rows that are
# already blank for the other "Resenha" columns G and J).
$newScores = @{
  2  = 0;  3  = 5;  4  = 7;  5  = $null; 6  = 10; 7  = 0;  8  = 7;
  9  = 5;  10 = 5;  11 = 5;  12 = 7;     13 = 5;  14 = $null; 15 = $null;
  16 = $null; 17 = 10; 18 = 5; 19 = 10;  20 = 0;  21 = 7;  22 = 10;
  23 = 0;  24 = 0;  25 = 0;  26 = 0;     27 = 10; 28 = 10; 29 = $null;
  30 = 0;  31 = 7;  32 = 0;  33 = 0;     34 = 7;  35 = 5;  36 = 10;
  37 = 0;  38 = 5;  39 = 7;  40 = 10;    41 = 0;  42 = 10; 43 = 7;  44 = 7
}

for ($row = 2; $row -le 44; $row++) {
    # Move the existing email address from L to the new M column.
    $email = $ws.Cells.Item($row, 12).Value2
    $ws.Cells.Item($row, 13).Value = $email

    # Write the new resenha score into L (blank where there is no score).
    $score = $newScores[$row]
    if ($null -eq $score) {
        $ws.Cells.Item($row, 12).Value = ""
    } else {
        $ws.Cells.Item($row, 12).Value = $score
    }
}

# Header row: L1 becomes the new resenha title, M1 keeps "Email".
$ws.Range("L1").Value = "Resenha Regime de Metas"
$ws.Range("M1").Value = "Email"
